# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The worker "YENDRYK ARNULFO PELUFFO RIVERA" (CC 1143389901) gets a new
# mora period (2011) inserted as the first data row, with an updated
# Valor Mora. The existing "KAREN MARGARITA SUAREZ GONZALEZ" (CC 1047497029)
# block is re-sorted into ascending period order (2305 .. 2408) below it,
# keeping each period's original Salario Basico / Valor Mora figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=16; Doc="1143389901"; Nombre="YENDRYK ARNULFO PELUFFO RIVERA"; Periodo="2011"; Salario=35960;  Mora=1754318},
    @{Row=17; Doc="1047497029"; Nombre="KAREN MARGARITA SUAREZ GONZALEZ"; Periodo="2305"; Salario=17013; Mora=1160000},
    @{Row=18; Doc="1047497029"; Nombre="KAREN MARGARITA SUAREZ GONZALEZ"; Periodo="2306"; Salario=46400; Mora=1160000},
    @{Row=19; Doc="1047497029"; Nombre="KAREN MARGARITA SUAREZ GONZALEZ"; Periodo="2307"; Salario=46400; Mora=1160000},
    @{Row=20; Doc="1047497029"; Nombre="KAREN MARGARITA SUAREZ GONZALEZ"; Periodo="2308"; Salario=46400; Mora=1160000},
    @{Row=21; Doc="1047497029"; Nombre="KAREN MARGARITA SUAREZ GONZALEZ"; Periodo="2309"; Salario=46400; Mora=1160000},
    @{Row=22; Doc="1047497029"; Nombre="KAREN MARGARITA SUAREZ GONZALEZ"; Periodo="2310"; Salario=46400; Mora=1160000},
    @{Row=23; Doc="1047497029"; Nombre="KAREN MARGARITA SUAREZ GONZALEZ"; Periodo="2311"; Salario=46400; Mora=1160000},
    @{Row=24; Doc="1047497029"; Nombre="KAREN MARGARITA SUAREZ GONZALEZ"; Periodo="2312"; Salario=46400; Mora=1160000},
    @{Row=25; Doc="1047497029"; Nombre="KAREN MARGARITA SUAREZ GONZALEZ"; Periodo="2401"; Salario=46400; Mora=1160000},
    @{Row=26; Doc="1047497029"; Nombre="KAREN MARGARITA SUAREZ GONZALEZ"; Periodo="2402"; Salario=46400; Mora=1160000},
    @{Row=27; Doc="1047497029"; Nombre="KAREN MARGARITA SUAREZ GONZALEZ"; Periodo="2403"; Salario=46400; Mora=1160000},
    @{Row=28; Doc="1047497029"; Nombre="KAREN MARGARITA SUAREZ GONZALEZ"; Periodo="2404"; Salario=46400; Mora=1160000},
    @{Row=29; Doc="1047497029"; Nombre="KAREN MARGARITA SUAREZ GONZALEZ"; Periodo="2405"; Salario=46400; Mora=1160000},
    @{Row=30; Doc="1047497029"; Nombre="KAREN MARGARITA SUAREZ GONZALEZ"; Periodo="2406"; Salario=46400; Mora=1160000},
    @{Row=31; Doc="1047497029"; Nombre="KAREN MARGARITA SUAREZ GONZALEZ"; Periodo="2407"; Salario=46400; Mora=1160000},
    @{Row=32; Doc="1047497029"; Nombre="KAREN MARGARITA SUAREZ GONZALEZ"; Periodo="2408"; Salario=40214; Mora=1160000}
)

foreach ($item in $data) {
    $r = $item.Row
    # Column B (Tipo Doc Trabajador = "CC") is unchanged for every row.
    $ws.Cells.Item($r, 3).Value = $item.Doc       # C - N° Doc Trabajador
    $ws.Cells.Item($r, 4).Value = $item.Nombre    # D - Nombre Trabajador
    $ws.Cells.Item($r, 5).Value = $item.Periodo   # E - Periodo Mora
    $ws.Cells.Item($r, 6).Value = $item.Salario   # F - Salario Basico
    $ws.Cells.Item($r, 7).Value = $item.Mora      # G - Valor Mora
}
